# "add project risks to register"
#
# Updates the "Risk Tracking Template" sheet (sheet 1):
#   - Retitle the template (A1)
#   - Extend the GDPR mitigation strategy note (G19)
#   - Add two new risk rows (17 & 18): "Key Team members sickness" and
#     "Scope Creep", matching the row-21/row-20 style pattern already used
#     in the table (red "High" fills on D/E/F, plain text elsewhere)
#   - Widen the Mitigation Strategy column (G) to fit the longer text
#   - Scroll the sheet view down so the new rows are visible

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk Tracking Template")

# --- Title -----------------------------------------------------------
$ws.Range("A1").Value = "Risk Tracking Template MAUIS"

# --- Extend existing GDPR mitigation text (row 19) --------------------
$ws.Cells.Item(19, 7).Value = "Design with privacy in mind, follow GDPR practices, Team Training"

# --- Row 22 (risk #17): Key Team members sickness ----------------------
$ws.Cells.Item(22, 2).Value = "Key Team members sickness"
$ws.Cells.Item(22, 3).Value = "Loss of team momentum "
$ws.Cells.Item(22, 4).Value = "High "
$ws.Cells.Item(22, 5).Value = "High "
$ws.Cells.Item(22, 6).Value = "High "
$ws.Cells.Item(22, 7).Value = "Design documentation for development, clear goals defined and backlog of tasks arranged into sprints"
$ws.Cells.Item(22, 8).Value = "Management"

# Reuse the existing "High" risk-level formatting (red fill) from row 6,
# column D, so D22:F22 pick up the same style already used elsewhere in
# the table instead of creating new duplicate styles.
$ws.Range("D6").Copy()
$ws.Range("D22:F22").PasteSpecial(-4122)

# --- Row 23 (risk #18): Scope Creep ------------------------------------
$ws.Cells.Item(23, 2).Value = "Scope Creep"
$ws.Cells.Item(23, 3).Value = "Goals become diluted with new less necessary tasks, original scope gets pushed back."
$ws.Cells.Item(23, 4).Value = "Medium"
$ws.Cells.Item(23, 5).Value = "High "
$ws.Cells.Item(23, 6).Value = "High "
$ws.Cells.Item(23, 7).Value = "Clear design documentation, team keeps to Agile principles and Scrum master and Product Owner keep development focussed and a tight feedback loop with stakeholders."
$ws.Cells.Item(23, 8).Value = "Scrum Master, ProductOwner"

# D23 is "Medium" risk level (yellow fill) - reuse that existing style too.
$ws.Range("E6").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("E23:F23").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Widen the Mitigation Strategy column to fit the longer text -------
$ws.Columns.Item(7).ColumnWidth = 80.83

# --- Scroll the view so the newly-added rows are visible ---------------
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D34").Select()
